$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.403.12"
$ws.Range("E2").Value = "  +0.86%  "

# Row 3
$ws.Range("D3").Value = "1.839.97"
$ws.Range("E3").Value = "  +3.47%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "225.15"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("D6").Value = "0.558"
$ws.Range("E6").Value = "  +1.71%  "

# Row 7
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "32.08"
$ws.Range("E8").Value = "  +1.51%  "

# Row 9
$ws.Range("E9").Value = "  +4.40%  "

# Row 10
$ws.Range("D10").Value = "0.0717"
$ws.Range("E10").Value = "  +9.33%  "

# Row 11
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  +0.47%  "

# Row 12
$ws.Range("D12").Value = "2.104.69"
$ws.Range("E12").Value = "  +3.46%  "

# Row 13
$ws.Range("D13").Value = "1.847.55"
$ws.Range("E13").Value = "  +4.03%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.650"
$ws.Range("E14").Value = "  +4.15%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "10.85"
$ws.Range("E15").Value = "  -1.85%  "

# Row 16
$ws.Range("D16").Value = "34.432.86"
$ws.Range("E16").Value = "  +0.90%  "

# Row 17
$ws.Range("D17").Value = "4.36"
$ws.Range("E17").Value = "  +3.87%  "

# Row 18
$ws.Range("D18").Value = "69.95"
$ws.Range("E18").Value = "  +1.80%  "

# Row 19
$ws.Range("D19").Value = "252.25"

# Row 20
$ws.Range("E20").Value = "  +8.34%  "

# Row 21
$ws.Range("D21").Value = "11.37"
$ws.Range("E21").Value = "  +9.99%  "

# Row 22
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("E23").Value = "  +2.63%  "

# Row 24
$ws.Range("E24").Value = "  +1.68%  "

# Row 25
$ws.Range("D25").Value = "160.91"
$ws.Range("E25").Value = "  +2.64%  "

# Row 26
$ws.Range("D26").Value = "16.73"
$ws.Range("E26").Value = "  +2.11%  "

# Row 27
$ws.Range("D27").Value = "7.28"
$ws.Range("E27").Value = "  +4.26%  "

# Row 28
$ws.Range("E28").Value = "  +2.25%  "

# Row 29
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("E30").Value = "  +5.18%  "

# Row 31
$ws.Range("D31").Value = "3.82"
$ws.Range("E31").Value = "  +1.58%  "

# Row 32
$ws.Range("E32").Value = "  +1.75%  "

# Row 33
$ws.Range("E33").Value = "  +1.38%  "

# Row 34
$ws.Range("E34").Value = "  +4.64%  "

# Row 35
$ws.Range("D35").Value = "1.459.65"
$ws.Range("E35").Value = "  +1.38%  "

# Row 36
$ws.Range("D36").Value = "0.647"
$ws.Range("E36").Value = "  +3.85%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0193"
$ws.Range("E37").Value = "  +3.35%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.06"
$ws.Range("E38").Value = "  +1.55%  "

# Row 39
$ws.Range("E39").Value = "  +9.38%  "

# Row 40
$ws.Range("D40").Value = "82.36"
$ws.Range("E40").Value = "  -0.56%  "

# Row 41
$ws.Range("E41").Value = "  -2.87%  "

# Row 42
$ws.Range("E42").Value = "  +0.34%  "

# Row 43
$ws.Range("E43").Value = "  +5.31%  "

# Row 44
$ws.Range("D44").Value = "6.09"
$ws.Range("E44").Value = "  +4.88%  "

# Row 45
$ws.Range("D45").Value = "2.002.63"
$ws.Range("E45").Value = "  +3.39%  "

# Row 46
$ws.Range("D46").Value = "0.0502"

# Row 47
$ws.Range("E47").Value = "  +0.93%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.50"
$ws.Range("E48").Value = "  +8.25%  "

# Row 49
$ws.Range("D49").Value = "12.08"
$ws.Range("E49").Value = "  +1.11%  "

# Row 50
$ws.Range("E50").Value = "  -0.08%  "

# Row 51
$ws.Range("E51").Value = "  +7.95%  "
